# This script applies a cyclic permutation of the data rows (rows 3-17) of the
# "Artfynd" worksheet: the *entire* content of each data row is relocated to a
# different row, while the header row (1) and the first data row (2) stay put.
#
# Because Excel's Range.Value2 array-assignment collapses "empty but present"
# string cells to $null (which then fails to round-trip as a cell at all), and
# because copying a full A:AY row range materializes blank cells for every
# column in that span (even ones that never had a cell to begin with), plain
# bulk copy/paste is not faithful enough here. Instead we:
#   1. snapshot every source row's *actual* cells (only the columns that are
#      really populated) into a staging area far below the used range,
#   2. wipe the destination rows clean,
#   3. copy, cell by cell, from the staging snapshot into the new row
#      position using Range.Copy (which preserves "empty string" cells
#      faithfully, unlike direct value assignment),
#   4. remove the staging rows again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are actually populated on every data row except row 10.
$standardCols = @("A","B","C","D","E","F","G","H","I","K","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AG","AT","AW","AX","AY")
# Row 10 additionally has J, N and AF populated.
$row10Cols = @("A","B","C","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AF","AG","AT","AW","AX","AY")

function Get-RowCols($rowNum) {
    if ($rowNum -eq 10) { return $row10Cols }
    return $standardCols
}

# Row permutation: destination row -> source row (values taken from the diff).
$mapping = @{
    3  = 14
    4  = 9
    5  = 10
    6  = 3
    7  = 11
    8  = 13
    9  = 6
    10 = 7
    11 = 16
    12 = 5
    13 = 17
    14 = 15
    15 = 8
    16 = 4
    17 = 12
}

$stagingOffset = 500

# Step 1: snapshot the current (original) content of rows 3-17 into staging rows.
foreach ($r in 3..17) {
    $stageRow = $r + $stagingOffset
    foreach ($c in (Get-RowCols $r)) {
        $src = $ws.Range("$c$r")
        $dst = $ws.Range("$c$stageRow")
        $src.Copy($dst)
    }
}

# Step 2: clear out the destination rows (every column that might currently hold data).
foreach ($r in 3..17) {
    $ws.Range("A" + $r + ":AY" + $r).ClearContents()
}

# Step 3: copy each row's snapshot into its new destination row per the mapping.
foreach ($destRow in 3..17) {
    $srcRow = $mapping[$destRow]
    $stageRow = $srcRow + $stagingOffset
    foreach ($c in (Get-RowCols $srcRow)) {
        $src = $ws.Range("$c$stageRow")
        $dst = $ws.Range("$c$destRow")
        $src.Copy($dst)
    }
}

# Step 4: remove the staging rows so they don't linger in the saved file.
foreach ($r in 3..17) {
    $stageRow = $r + $stagingOffset
    $ws.Range("A" + $stageRow + ":AY" + $stageRow).ClearContents()
}
